$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns B-E, values + same style as A1
$ws.Range("B1").Value = "Am 03"
$ws.Range("C1").Value = "Am 04"
$ws.Range("D1").Value = "Am 07"
$ws.Range("E1").Value = "Am 08"

$ws.Range("A1").Copy()
$ws.Range("B1:E1").PasteSpecial(-4122)  # xlPasteFormats

# Row 2 - Peso Inicial (g)
$ws.Range("B2").Value = 0.777
$ws.Range("C2").Value = 0.768
$ws.Range("D2").Value = 0.77
$ws.Range("E2").Value = 0.78

# Row 3 - Fase (0d - 21d)
$ws.Range("B3").Value = "+21dias"
$ws.Range("C3").Value = "+21dias"
$ws.Range("D3").Value = "+21dias"
$ws.Range("E3").Value = "+21dias"

# Row 4 - Peso (g)
$ws.Range("B4").Value = 0.77
$ws.Range("C4").Value = 0.762
$ws.Range("D4").Value = 0.765
$ws.Range("E4").Value = 0.776

# Row 5 - Wc 21d [g/Ah]
$ws.Range("B5").Value = 1.483
$ws.Range("C5").Value = 1.218
$ws.Range("D5").Value = 1.002
$ws.Range("E5").Value = 0.802
